$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 0. Stash the old summary row's cell style (row 5, style s="5") in a
#    scratch cell on row 1 (an untouched row, so later row-shifting
#    operations on rows >= 5 won't move it) for reuse by the new
#    summary row (row 7) further down.
# -----------------------------------------------------------------------
$ws.Range("F5").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# -----------------------------------------------------------------------
# 1. Remove the old "Все классы - 1 учеников" summary row (row 5). It is
#    currently the last row in the sheet, so this simply shrinks the
#    sheet rather than disturbing any other data.
# -----------------------------------------------------------------------
$ws.Rows(5).Delete()

# -----------------------------------------------------------------------
# 2. Update row 2 (existing student's data changed).
#    Several cells switch from empty/numeric to free-form text; a leading
#    apostrophe forces text storage, then the whole row's formatting is
#    re-applied from a pristine style-4 cell (D2, itself about to be
#    overwritten) so the stray "quote prefix" flag added by the
#    apostrophe doesn't create a new/different style.
# -----------------------------------------------------------------------
$ws.Range("B2").Value = "erwyg dfgs dfgs"
$ws.Range("C2").Value = "'856477567657"
$ws.Range("D2").Value = "'dfsa ddafs asdf"
$ws.Range("E2").Value = "'+7 (132) 412-34-31"
$ws.Range("F2").Value = "'sadf"
$ws.Range("G2").Value = 3124
$ws.Range("H2").Value = "3 Б"

$ws.Range("H2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# -----------------------------------------------------------------------
# 3. Add row 3 (new student), formatting copied from row 2.
# -----------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Dolmagambetov Karen "
$ws.Range("C3").Value = "'010608550491"
$ws.Range("D3").Value = "Долмагамбетов Талгат Аманжолович"
$ws.Range("E3").Value = "'+7 (771) 168-86-87"
$ws.Range("F3").Value = "'214234"
$ws.Range("G3").Value = 213
$ws.Range("H3").Value = "4 В"

$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

# -----------------------------------------------------------------------
# 4. Add row 4 (new student), formatting copied from row 2.
# -----------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Сапаров Айбек Галымжанович"
$ws.Range("C4").Value = "'010101854353"
$ws.Range("D4").Value = "Сапаров Галымжан Талгатович"
$ws.Range("E4").Value = "'+7 (771) 168-86-87"
$ws.Range("F4").Value = "Ораза Татеулы 13А"
$ws.Range("G4").Value = 80000
$ws.Range("H4").Value = "6 А"

$ws.Range("A2:H2").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)

# -----------------------------------------------------------------------
# 5. Add the new summary row (row 7), reusing the style stashed in step 0.
# -----------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = "Все классы - 3 учеников"

# -----------------------------------------------------------------------
# 6. Clean up the scratch cell so it leaves no trace in the saved sheet.
# -----------------------------------------------------------------------
$ws.Range("Z1").ClearFormats()
$ws.Range("Z1").Value = "x"
$ws.Range("Z1").Value = ""

$ws.Range("A1").Select()
